$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 10
$ws.Range("N2").Value = 2.01
$ws.Range("O2").Value = 1.89

# Row 5
$ws.Range("G5").Value = 2.6
$ws.Range("H5").Value = 2.7
$ws.Range("I5").Value = 3.2
$ws.Range("T5").Value = 6
$ws.Range("U5").Value = 11
$ws.Range("V5").Value = 12
$ws.Range("Y5").Value = 51
$ws.Range("AA5").Value = 5.5
$ws.Range("AE5").Value = 6.5
$ws.Range("AF5").Value = 13

# Row 21
$ws.Range("G21").Value = 3.5
$ws.Range("H21").Value = 3.8
$ws.Range("I21").Value = 1.95
$ws.Range("R21").Value = 1.53
$ws.Range("S21").Value = 2.38
$ws.Range("U21").Value = 21
$ws.Range("V21").Value = 12
$ws.Range("X21").Value = 26
$ws.Range("Y21").Value = 26
$ws.Range("AA21").Value = 7.5
$ws.Range("AB21").Value = 12
$ws.Range("AC21").Value = 34
$ws.Range("AF21").Value = 11
$ws.Range("AH21").Value = 19

# Row 32
$ws.Range("G32").Value = 3.25
$ws.Range("H32").Value = 2.8
$ws.Range("I32").Value = 2.42
$ws.Range("Q32").Value = 2.65
$ws.Range("Y32").Value = 37
$ws.Range("AE32").Value = 7.4
$ws.Range("AF32").Value = 12
$ws.Range("AH32").Value = 27

# Row 34
$ws.Range("J34").Value = 1.08
$ws.Range("K34").Value = 8

# Row 37
$ws.Range("G37").Value = 2.7
$ws.Range("H37").Value = 3.9
$ws.Range("I37").Value = 2.27
$ws.Range("L37").Value = 1.16
$ws.Range("M37").Value = 4.8
$ws.Range("N37").Value = 1.52
$ws.Range("O37").Value = 2.42
$ws.Range("P37").Value = 1.26
$ws.Range("Q37").Value = 3.5
$ws.Range("R37").Value = 1.48
$ws.Range("S37").Value = 2.55
$ws.Range("T37").Value = 11
$ws.Range("U37").Value = 14
$ws.Range("V37").Value = 8.6
$ws.Range("W37").Value = 28
$ws.Range("X37").Value = 16
$ws.Range("Y37").Value = 19
$ws.Range("Z37").Value = 15
$ws.Range("AA37").Value = 6.6
$ws.Range("AB37").Value = 9.8
$ws.Range("AC37").Value = 35
$ws.Range("AD37").Value = 101
$ws.Range("AE37").Value = 10
$ws.Range("AF37").Value = 12
$ws.Range("AG37").Value = 7.8
$ws.Range("AH37").Value = 21
$ws.Range("AI37").Value = 14
$ws.Range("AJ37").Value = 18

# Row 43
$ws.Range("G43").Value = 2.55
$ws.Range("I43").Value = 2.8
$ws.Range("K43").Value = 7.5
$ws.Range("W43").Value = 26

# Row 46
$ws.Range("K46").Value = 8
$ws.Range("N46").Value = 2.25
$ws.Range("O46").Value = 1.62

# Row 53
$ws.Range("N53").Value = 1.88
$ws.Range("O53").Value = 1.98

# Row 69
$ws.Range("G69").Value = 1.95
$ws.Range("I69").Value = 3.75
$ws.Range("J69").Value = 1.06
$ws.Range("K69").Value = 10
$ws.Range("N69").Value = 2.05
$ws.Range("O69").Value = 1.8
$ws.Range("T69").Value = 7
$ws.Range("V69").Value = 9
$ws.Range("X69").Value = 17
$ws.Range("Y69").Value = 29
$ws.Range("AI69").Value = 29
